$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds years 2000-2020 in rows 2-21.
# Target: drop the 2000-2009 rows (old rows 2-11) so 2010-2020 shift up
# to rows 2-11, then append a new 2021 row as row 12.

# Delete old rows 2-11 (2000年..2009年) -- this shifts rows 12-21 up to 2-11.
$ws.Rows("2:11").Delete()

# Give the new row 12 the same formatting as the row above it (row 11,
# which now holds the 2020年 data) before filling in the 2021年 figures.
$ws.Range("A11").Copy($ws.Range("A12"))

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 715
$ws.Range("C12").Value = 89
$ws.Range("D12").Value = 77
$ws.Range("E12").Value = 38
$ws.Range("F12").Value = 120
$ws.Range("G12").Value = 30
